# Frontend_Table_Description.xlsx edit:
# Insert a new "record_id_42076f" column-description row just above the
# existing row 13 (the separator row that precedes the "fall" table block),
# pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("frontend_table_description")

# Insert a new, blank row at row 13 - everything currently at row 13 and
# below (the "fall" table section, etc.) shifts down by one row.
$ws.Range("A13").EntireRow.Insert()

# Populate the newly inserted row 13 with the temporary RedCap column.
$ws.Range("B13").Value = "record_id_42076f"
$ws.Range("C13").Value = "ToDo: zu Entfernen wenn aus RC Projekt"
$ws.Range("D13").Value = "varchar"

# Restore the view state: scrolled so row 11 is at the top, with C22
# (now the "fall_studienphase" description cell) selected.
$window = $ws.Application.ActiveWindow
$window.ScrollRow = 11
$ws.Range("C22").Select()
